# Update the cryptocurrency price/volume table with the latest scraped
# figures (GitHub Actions refresh). Row 48/49 also swap contents because the
# source ranking re-ordered "Maker" and "ordi".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Col = 4; Value = '43.974.40'; ForceText = $false },
    @{ Row = 2; Col = 5; Value = '  +4.23%  '; ForceText = $false },
    @{ Row = 3; Col = 4; Value = '2.349.19'; ForceText = $false },
    @{ Row = 3; Col = 5; Value = '  +3.11%  '; ForceText = $false },
    @{ Row = 4; Col = 5; Value = '  -0.83%  '; ForceText = $false },
    @{ Row = 5; Col = 4; Value = '314.28'; ForceText = $true },
    @{ Row = 5; Col = 5; Value = '  +0.62%  '; ForceText = $false },
    @{ Row = 6; Col = 4; Value = '109.42'; ForceText = $true },
    @{ Row = 6; Col = 5; Value = '  +8.01%  '; ForceText = $false },
    @{ Row = 7; Col = 4; Value = '0.635'; ForceText = $true },
    @{ Row = 7; Col = 5; Value = '  +3.67%  '; ForceText = $false },
    @{ Row = 8; Col = 5; Value = '  -0.03%  '; ForceText = $false },
    @{ Row = 9; Col = 4; Value = '0.625'; ForceText = $true },
    @{ Row = 9; Col = 5; Value = '  +5.78%  '; ForceText = $false },
    @{ Row = 10; Col = 4; Value = '41.74'; ForceText = $true },
    @{ Row = 10; Col = 5; Value = '  +8.34%  '; ForceText = $false },
    @{ Row = 11; Col = 4; Value = '0.0923'; ForceText = $true },
    @{ Row = 11; Col = 5; Value = '  +3.01%  '; ForceText = $false },
    @{ Row = 12; Col = 4; Value = '8.62'; ForceText = $true },
    @{ Row = 12; Col = 5; Value = '  +5.06%  '; ForceText = $false },
    @{ Row = 13; Col = 4; Value = '1.01'; ForceText = $true },
    @{ Row = 13; Col = 5; Value = '  +4.08%  '; ForceText = $false },
    @{ Row = 14; Col = 5; Value = '  -0.03%  '; ForceText = $false },
    @{ Row = 15; Col = 4; Value = '15.54'; ForceText = $true },
    @{ Row = 15; Col = 5; Value = '  +3.80%  '; ForceText = $false },
    @{ Row = 16; Col = 4; Value = '2.695.13'; ForceText = $false },
    @{ Row = 16; Col = 5; Value = '  +2.71%  '; ForceText = $false },
    @{ Row = 17; Col = 4; Value = '2.348.23'; ForceText = $false },
    @{ Row = 17; Col = 5; Value = '  +1.03%  '; ForceText = $false },
    @{ Row = 18; Col = 4; Value = '44.113.24'; ForceText = $false },
    @{ Row = 18; Col = 5; Value = '  +4.23%  '; ForceText = $false },
    @{ Row = 19; Col = 4; Value = '7.61'; ForceText = $true },
    @{ Row = 19; Col = 5; Value = '  +5.55%  '; ForceText = $false },
    @{ Row = 20; Col = 5; Value = '  +2.96%  '; ForceText = $false },
    @{ Row = 21; Col = 4; Value = '13.00'; ForceText = $true },
    @{ Row = 21; Col = 5; Value = '  -3.35%  '; ForceText = $false },
    @{ Row = 22; Col = 4; Value = '74.69'; ForceText = $true },
    @{ Row = 22; Col = 5; Value = '  +2.61%  '; ForceText = $false },
    @{ Row = 23; Col = 5; Value = '  +0.60%  '; ForceText = $false },
    @{ Row = 24; Col = 4; Value = '268.91'; ForceText = $true },
    @{ Row = 24; Col = 5; Value = '  +2.10%  '; ForceText = $false },
    @{ Row = 25; Col = 4; Value = '2.29'; ForceText = $true },
    @{ Row = 25; Col = 5; Value = '  +6.34%  '; ForceText = $false },
    @{ Row = 27; Col = 4; Value = '7.60'; ForceText = $true },
    @{ Row = 27; Col = 5; Value = '  +11.24%  '; ForceText = $false },
    @{ Row = 28; Col = 4; Value = '11.19'; ForceText = $true },
    @{ Row = 28; Col = 5; Value = '  +5.63%  '; ForceText = $false },
    @{ Row = 29; Col = 4; Value = '2.32'; ForceText = $true },
    @{ Row = 29; Col = 5; Value = '  -0.76%  '; ForceText = $false },
    @{ Row = 30; Col = 4; Value = '39.79'; ForceText = $true },
    @{ Row = 30; Col = 5; Value = '  +11.31%  '; ForceText = $false },
    @{ Row = 31; Col = 4; Value = '22.72'; ForceText = $true },
    @{ Row = 31; Col = 5; Value = '  +2.07%  '; ForceText = $false },
    @{ Row = 32; Col = 4; Value = '168.95'; ForceText = $true },
    @{ Row = 32; Col = 5; Value = '  +2.61%  '; ForceText = $false },
    @{ Row = 33; Col = 4; Value = '0.0913'; ForceText = $true },
    @{ Row = 33; Col = 5; Value = '  +6.53%  '; ForceText = $false },
    @{ Row = 34; Col = 5; Value = '  +9.56%  '; ForceText = $false },
    @{ Row = 35; Col = 4; Value = '0.133'; ForceText = $true },
    @{ Row = 35; Col = 5; Value = '  +2.39%  '; ForceText = $false },
    @{ Row = 36; Col = 5; Value = '  +4.46%  '; ForceText = $false },
    @{ Row = 37; Col = 4; Value = '4.75'; ForceText = $true },
    @{ Row = 37; Col = 5; Value = '  +6.96%  '; ForceText = $false },
    @{ Row = 38; Col = 4; Value = '0.0366'; ForceText = $true },
    @{ Row = 38; Col = 5; Value = '  +5.64%  '; ForceText = $false },
    @{ Row = 39; Col = 5; Value = '  +10.99%  '; ForceText = $false },
    @{ Row = 40; Col = 5; Value = '  +2.57%  '; ForceText = $false },
    @{ Row = 41; Col = 4; Value = '1.74'; ForceText = $true },
    @{ Row = 41; Col = 5; Value = '  +11.74%  '; ForceText = $false },
    @{ Row = 42; Col = 4; Value = '105.32'; ForceText = $true },
    @{ Row = 42; Col = 5; Value = '  +7.42%  '; ForceText = $false },
    @{ Row = 43; Col = 4; Value = '13.83'; ForceText = $true },
    @{ Row = 43; Col = 5; Value = '  +16.66%  '; ForceText = $false },
    @{ Row = 44; Col = 4; Value = '0.241'; ForceText = $true },
    @{ Row = 44; Col = 5; Value = '  +7.25%  '; ForceText = $false },
    @{ Row = 45; Col = 4; Value = '71.99'; ForceText = $true },
    @{ Row = 45; Col = 5; Value = '  +4.89%  '; ForceText = $false },
    @{ Row = 46; Col = 5; Value = '  -0.06%  '; ForceText = $false },
    @{ Row = 47; Col = 4; Value = '115.67'; ForceText = $true },
    @{ Row = 47; Col = 5; Value = '  +5.37%  '; ForceText = $false },
    @{ Row = 48; Col = 2; Value = 'ordi'; ForceText = $false },
    @{ Row = 48; Col = 3; Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'; ForceText = $false },
    @{ Row = 48; Col = 4; Value = '77.55'; ForceText = $true },
    @{ Row = 48; Col = 5; Value = '  -1.72%  '; ForceText = $false },
    @{ Row = 49; Col = 2; Value = 'Maker'; ForceText = $false },
    @{ Row = 49; Col = 3; Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; ForceText = $false },
    @{ Row = 49; Col = 4; Value = '1.658.75'; ForceText = $false },
    @{ Row = 49; Col = 5; Value = '  -2.63%  '; ForceText = $false },
    @{ Row = 50; Col = 5; Value = '  +4.70%  '; ForceText = $false },
    @{ Row = 51; Col = 5; Value = '  +17.90%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Cells.Item($u.Row, $u.Col)
    if ($u.ForceText) {
        # Column D sometimes holds values that look numeric (e.g. "314.28").
        # Force a text format before/after the write so the engine keeps the
        # exact literal instead of silently coercing it into a Double (which
        # would also mangle values like "0.0920" -> "0.092").
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
